# Issue #35 Blank settings page
# Applies the Issues log updates: re-numbers a couple of "Requires" links,
# marks several completed/duplicate issues as DONE & hidden, rewrites the
# note on issue #13 (Better Top Navigation), extends the AutoFilter /
# _FilterDatabase range, and appends a new issue #35 "Blank settings page".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# --- Turn off the existing AutoFilter before touching data, so that when it
#     is re-applied below it does not silently swallow the brand-new row 35.
$ws.AutoFilterMode = $false

# --- Row 3 (issue #2 "Speed"): Requires now points at issue #35 not #24,
#     and Priority moves from 2 to 1.
$ws.Range("B3").Value = 1
$ws.Range("F3").Value = 35

# --- Row 4 (issue #3 "Brightness"): same Requires / Priority update, and it
#     now also gets a Requires value (previously blank).
$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 35

# --- Row 14 (issue #13): now marked DONE & hidden by the filter, its note
#     is expanded, and its row height is forced to 29 (2 wrapped lines).
$ws.Range("C14").Value = "DONE"
$ws.Range("H14").Value = "Better Top Navigation " + [char]0x2026 + " title and hamburger." + "`n" + "No code change " + [char]0x2026 + " fixed as consequnce of other changes"
$ws.Rows.Item(14).RowHeight = 29
$ws.Rows.Item(14).EntireRow.Hidden = $true

# --- Row 27 (issue #27 "Hearbeat"): Requires now points at issue #35.
$ws.Range("F27").Value = 35

# --- Rows 28, 29, 30: existing DONE issues that are now also hidden by the
#     filter (no content changes).
$ws.Rows.Item(28).EntireRow.Hidden = $true
$ws.Rows.Item(29).EntireRow.Hidden = $true
$ws.Rows.Item(30).EntireRow.Hidden = $true

# --- Row 34 (issue #33 "icon config"): marked DONE and hidden.
$ws.Range("C34").Value = "DONE"
$ws.Rows.Item(34).EntireRow.Hidden = $true

# --- New row 35: issue #35 "Blank settings page".
$ws.Range("A35").Value = 35
$ws.Range("B35").Value = 1
$ws.Range("C35").Value = "DONE"
$ws.Range("H35").Value = "Settings page placehold for other stuff"
$ws.Range("E35").Value = "Blank settings page"

# --- Re-apply the AutoFilter over the (still 34-row) table range, restoring
#     the "(Blanks)" filter on the Status column (colId 2 => column C).
$ws.Range("A1:H34").AutoFilter(3, @(""), 7)

# --- Update the hidden _FilterDatabase defined name to match the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Issues!_FilterDatabase") {
        $n.RefersTo = "=Issues!`$A`$1:`$H`$34"
    }
}

# --- Move the active selection to E37, matching the saved view state.
$ws.Activate()
$ws.Range("E37").Select()
